# chore: update Sheets via scheduled runner
# Refreshes the currentAveragePrice / LevePrice / LeveProfit columns
# (H:N) for a handful of leve rows across all eight Tonberry_Profits
# sheets, reflecting newer market-board pricing data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1523.7142
$ws.Range("J17").Value = 1292.9181
$ws.Range("L17").Value = 3878.754300000001
$ws.Range("N17").Value = -4214.754300000001

$ws.Range("H40").Value = 1024.3778
$ws.Range("I40").Value = 1005.5897
$ws.Range("J40").Value = 1146.5
$ws.Range("K40").Value = 1005.5897
$ws.Range("L40").Value = 1146.5
$ws.Range("M40").Value = -830.5897
$ws.Range("N40").Value = -1496.5

$ws.Range("H113").Value = 142958430
$ws.Range("I113").Value = 118167.5
$ws.Range("K113").Value = 118167.5
$ws.Range("M113").Value = -114913.5

$ws.Range("H116").Value = 4528.4287
$ws.Range("I116").Value = 2566.3333
$ws.Range("K116").Value = 2566.3333
$ws.Range("M116").Value = 875.6667000000002

$ws.Range("H127").Value = 2089.3684
$ws.Range("I127").Value = 1764.4286
$ws.Range("K127").Value = 5293.2858
$ws.Range("M127").Value = -333.2857999999997

$ws.Range("H138").Value = 1592.5302
$ws.Range("I138").Value = 1283.6666
$ws.Range("J138").Value = 2398.261
$ws.Range("K138").Value = 3850.9998
$ws.Range("L138").Value = 7194.782999999999
$ws.Range("M138").Value = 1289.0002
$ws.Range("N138").Value = -17474.783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5555555
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 3761.736
$ws.Range("I32").Value = 3065.111
$ws.Range("K32").Value = 3065.111
$ws.Range("M32").Value = -2778.111

$ws.Range("H110").Value = 1775.6
$ws.Range("I110").Value = 1548.3077
$ws.Range("K110").Value = 1548.3077
$ws.Range("M110").Value = 496.6922999999999

$ws.Range("H116").Value = 5555555
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 1491.5135
$ws.Range("I132").Value = 1249.8148
$ws.Range("J132").Value = 2144.1
$ws.Range("K132").Value = 3749.4444
$ws.Range("L132").Value = 6432.299999999999
$ws.Range("M132").Value = -1219.4444
$ws.Range("N132").Value = -11492.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5555555
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H134").Value = 1353.0625
$ws.Range("I134").Value = 1353.0625
$ws.Range("K134").Value = 4059.1875
$ws.Range("M134").Value = -1524.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1055.2106
$ws.Range("J22").Value = 1727.7778
$ws.Range("L22").Value = 1727.7778
$ws.Range("N22").Value = -2427.7778

$ws.Range("H31").Value = 2233818.8
$ws.Range("I31").Value = 3572660.5
$ws.Range("J31").Value = 2415.8333
$ws.Range("K31").Value = 3572660.5
$ws.Range("L31").Value = 2415.8333
$ws.Range("M31").Value = -3572365.5
$ws.Range("N31").Value = -3005.8333

$ws.Range("H34").Value = 2233818.8
$ws.Range("I34").Value = 3572660.5
$ws.Range("J34").Value = 2415.8333
$ws.Range("K34").Value = 3572660.5
$ws.Range("L34").Value = 2415.8333
$ws.Range("M34").Value = -3572458.5
$ws.Range("N34").Value = -2819.8333

$ws.Range("H58").Value = 1978419.5
$ws.Range("I58").Value = 2899429
$ws.Range("K58").Value = 2899429
$ws.Range("M58").Value = -2899226

$ws.Range("H132").Value = 1439.4524
$ws.Range("I132").Value = 986.76666
$ws.Range("J132").Value = 2571.1667
$ws.Range("K132").Value = 2960.29998
$ws.Range("L132").Value = 7713.500100000001
$ws.Range("M132").Value = -430.2999799999998
$ws.Range("N132").Value = -12773.5001

$ws.Range("H134").Value = 1430.678
$ws.Range("J134").Value = 2026.4166
$ws.Range("L134").Value = 6079.2498
$ws.Range("N134").Value = -11149.2498

$ws.Range("H136").Value = 1978419.5
$ws.Range("I136").Value = 2899429
$ws.Range("K136").Value = 8698287
$ws.Range("M136").Value = -8695737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I6").Value = 462.42856
$ws.Range("K6").Value = 1387.28568
$ws.Range("M6").Value = -1274.28568

$ws.Range("H33").Value = 60.846153
$ws.Range("I33").Value = 46.833332
$ws.Range("J33").Value = 229
$ws.Range("K33").Value = 280.999992
$ws.Range("L33").Value = 1374
$ws.Range("M33").Value = 2.00000799999998
$ws.Range("N33").Value = -1940

$ws.Range("H76").Value = 3506.5
$ws.Range("I76").Value = 2013
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 6039
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -5656
$ws.Range("N76").Value = -15766

$ws.Range("H79").Value = 3506.5
$ws.Range("I79").Value = 2013
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 6039
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -4713
$ws.Range("N79").Value = -17652

$ws.Range("H107").Value = 386.25
$ws.Range("J107").Value = 412
$ws.Range("L107").Value = 1236
$ws.Range("N107").Value = -5076

$ws.Range("H131").Value = 6955010
$ws.Range("I131").Value = 35714740
$ws.Range("J131").Value = 13006.482
$ws.Range("K131").Value = 107144220
$ws.Range("L131").Value = 39019.446
$ws.Range("M131").Value = -107139180
$ws.Range("N131").Value = -49099.446

$ws.Range("H132").Value = 1059.6
$ws.Range("J132").Value = 1449
$ws.Range("L132").Value = 13041
$ws.Range("N132").Value = -18101

$ws.Range("H140").Value = 2142.532
$ws.Range("I140").Value = 1363.3684
$ws.Range("J140").Value = 2671.25
$ws.Range("K140").Value = 4090.1052
$ws.Range("L140").Value = 8013.75
$ws.Range("M140").Value = 1089.8948
$ws.Range("N140").Value = -18373.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4299.1665
$ws.Range("J70").Value = 4299.1665
$ws.Range("L70").Value = 4299.1665
$ws.Range("N70").Value = -4839.1665

$ws.Range("H73").Value = 4299.1665
$ws.Range("J73").Value = 4299.1665
$ws.Range("L73").Value = 4299.1665
$ws.Range("N73").Value = -6171.1665

$ws.Range("H113").Value = 1370.75
$ws.Range("I113").Value = 1435.5714
$ws.Range("J113").Value = 1280
$ws.Range("K113").Value = 1435.5714
$ws.Range("L113").Value = 1280
$ws.Range("M113").Value = 734.4286
$ws.Range("N113").Value = -5620

$ws.Range("H122").Value = 1195.125
$ws.Range("J122").Value = 1998
$ws.Range("L122").Value = 5994
$ws.Range("N122").Value = -10894

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2352.1177
$ws.Range("I61").Value = 1999.0667
$ws.Range("K61").Value = 1999.0667
$ws.Range("M61").Value = -1797.0667

$ws.Range("H82").Value = 2139.4167
$ws.Range("J82").Value = 2638.4285
$ws.Range("L82").Value = 2638.4285
$ws.Range("N82").Value = -3360.4285

$ws.Range("H85").Value = 2139.4167
$ws.Range("J85").Value = 2638.4285
$ws.Range("L85").Value = 2638.4285
$ws.Range("N85").Value = -5134.4285

$ws.Range("H113").Value = 2352.1177
$ws.Range("I113").Value = 1999.0667
$ws.Range("K113").Value = 1999.0667
$ws.Range("M113").Value = 170.9332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5327.2593
$ws.Range("J14").Value = 4949.5
$ws.Range("L14").Value = 4949.5
$ws.Range("N14").Value = -5285.5

$ws.Range("H95").Value = 99992.2
$ws.Range("J95").Value = 99992.2
$ws.Range("L95").Value = 99992.2
$ws.Range("N95").Value = -105484.2

$ws.Range("H126").Value = 5961.1113
$ws.Range("I126").Value = 7538.125
$ws.Range("J126").Value = 3667.2727
$ws.Range("K126").Value = 22614.375
$ws.Range("L126").Value = 11001.8181
$ws.Range("M126").Value = -20144.375
$ws.Range("N126").Value = -15941.8181

$ws.Range("H132").Value = 1217.9841
$ws.Range("I132").Value = 949.9423
$ws.Range("K132").Value = 2849.8269
$ws.Range("M132").Value = -319.8269
